# Refitting NCDEs to individual patients (for manuscript figure)
# Adds a "Label" column (H) to the classification results sheet:
#   Label = 0 for Control patients, 1 for MDD patients

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1, formatted like the other header cells (copy format from G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Label"

# New data column H2:H21 — 0 for Control rows, 1 for MDD rows
# (rows 2-11 and 12-21 each cycle through the five Control then five MDD patients)
$labels = @(0, 0, 0, 0, 0, 1, 1, 1, 1, 1, 0, 0, 0, 0, 0, 1, 1, 1, 1, 1)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $labels[$i]
}
